$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.083.37"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.44%  '
$ws.Range("D3").Value = "'3.295.01"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.65%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").Value = "'587.17"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.59%  '
$ws.Range("D6").Value = "'179.19"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.16%  '
$ws.Range("D7").Value = "'0.641"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.95%  '
$ws.Range("E8").Value = '  +0.10%  '
$ws.Range("D9").Value = "'3.295.05"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.72%  '
$ws.Range("E10").Value = '  -0.43%  '
$ws.Range("D11").Value = "'6.84"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.08%  '
$ws.Range("D12").Value = "'0.401"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.70%  '
$ws.Range("D13").Value = "'3.872.39"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.69%  '
$ws.Range("E14").Value = '  -2.51%  '
$ws.Range("D15").Value = "'66.205.22"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.50%  '
$ws.Range("D16").Value = "'26.54"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.46%  '
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = "'3.312.29"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.49%  '
$ws.Range("B18").Value = 'ShibaInu'
$ws.Range("C18").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D18").Value = "'0.0000163"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.22%  '
$ws.Range("D19").Value = "'423.30"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.94%  '
$ws.Range("D20").Value = "'5.48"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.35%  '
$ws.Range("D21").Value = "'13.03"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.83%  '
$ws.Range("D22").Value = "'7.28"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.38%  '
$ws.Range("D23").Value = "'0.997"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.28%  '
$ws.Range("D24").Value = "'71.26"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.49%  '
$ws.Range("E25").Value = '  -0.05%  '
$ws.Range("B26").Value = 'Kaspa'
$ws.Range("C26").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D26").Value = "'0.207"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +6.20%  '
$ws.Range("B27").Value = 'Polygon'
$ws.Range("C27").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D27").Value = "'0.509"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.64%  '
$ws.Range("E28").Value = '  +1.34%  '
$ws.Range("D29").Value = "'9.46"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +6.72%  '
$ws.Range("E30").Value = '  -0.05%  '
$ws.Range("E31").Value = '  -0.45%  '
$ws.Range("D32").Value = "'22.27"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.18%  '
$ws.Range("E33").Value = '  +0.04%  '
$ws.Range("D34").Value = "'5.15"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.54%  '
$ws.Range("D35").Value = "'6.58"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.03%  '
$ws.Range("E36").Value = '  +0.61%  '
$ws.Range("D37").Value = "'157.63"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.14%  '
$ws.Range("E38").Value = '  -0.07%  '
$ws.Range("D39").Value = "'2.866.48"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.40%  '
$ws.Range("D40").Value = "'1.79"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.67%  '
$ws.Range("D41").Value = "'26.32"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.52%  '
$ws.Range("D42").Value = "'4.33"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.56%  '
$ws.Range("D43").Value = "'0.749"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.13%  '
$ws.Range("D44").Value = "'39.64"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.56%  '
$ws.Range("D45").Value = "'5.90"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.27%  '
$ws.Range("D46").Value = "'2.30"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.76%  '
$ws.Range("D47").Value = "'0.0638"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.46%  '
$ws.Range("D48").Value = "'312.75"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.62%  '
$ws.Range("D49").Value = "'22.89"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.87%  '
$ws.Range("D50").Value = "'0.0269"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.76%  '
$ws.Range("E51").Value = '  -0.15%  '
